# Added code for selenium grid setup
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 identity/pass values
$ws.Range("A2").Value = "sumitIdentity7"
$ws.Range("B2").Value = "SummitPass7"

# Add new row 3 identity/pass values
$ws.Range("A3").Value = "sumitIdentity8"
$ws.Range("B3").Value = "SummitPass8"

# Move the selection to E2
$null = $ws.Range("E2").Select()
